# Applies the cryptos.xlsx price/volume refresh described in the commit.
# D (Price) values are assigned with a forced leading apostrophe so Excel
# keeps them as text (matching the original inlineStr cells) instead of
# normalizing numeric-looking strings such as "1.00" -> 1 or "18.40" -> 18.4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''51.957.52'
$ws.Range("E2").Value = '  +0.98%  '
# Row 3
$ws.Range("D3").Value = '''2.823.74'
$ws.Range("E3").Value = '  +2.69%  '
# Row 4
$ws.Range("E4").Value = '  +0.03%  '
# Row 5
$ws.Range("D5").Value = '''355.75'
$ws.Range("E5").Value = '  +6.98%  '
# Row 6
$ws.Range("D6").Value = '''113.52'
$ws.Range("E6").Value = '  -1.81%  '
# Row 7
$ws.Range("D7").Value = '''0.553'
$ws.Range("E7").Value = '  +3.02%  '
# Row 8
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.03%  '
# Row 9
$ws.Range("E9").Value = '  +4.72%  '
# Row 10
$ws.Range("D10").Value = '''41.88'
$ws.Range("E10").Value = '  +0.86%  '
# Row 11
$ws.Range("E11").Value = '  -0.37%  '
# Row 12
$ws.Range("D12").Value = '''20.05'
$ws.Range("E12").Value = '  -0.61%  '
# Row 13
$ws.Range("E13").Value = '  +1.47%  '
# Row 14
$ws.Range("D14").Value = '''7.73'
# Row 15
$ws.Range("D15").Value = '''3.245.59'
$ws.Range("E15").Value = '  +2.09%  '
# Row 16
$ws.Range("D16").Value = '''2.829.53'
$ws.Range("E16").Value = '  +3.44%  '
# Row 17
$ws.Range("D17").Value = '''0.899'
$ws.Range("E17").Value = '  +1.85%  '
# Row 18
$ws.Range("D18").Value = '''51.895.85'
$ws.Range("E18").Value = '  +0.84%  '
# Row 19
$ws.Range("D19").Value = '''7.42'
$ws.Range("E19").Value = '  +8.44%  '
# Row 20
$ws.Range("E20").Value = '  -1.62%  '
# Row 21
$ws.Range("D21").Value = '''13.61'
$ws.Range("E21").Value = '  +1.47%  '
# Row 22
$ws.Range("D22").Value = '''0.0₃0996'
$ws.Range("E22").Value = '  +2.43%  '
# Row 23
$ws.Range("D23").Value = '''270.33'
$ws.Range("E23").Value = '  -2.81%  '
# Row 24
$ws.Range("D24").Value = '''69.79'
$ws.Range("E24").Value = '  +0.51%  '
# Row 25
$ws.Range("D25").Value = '''2.79'
$ws.Range("E25").Value = '  +5.29%  '
# Row 26
$ws.Range("D26").Value = '''26.82'
$ws.Range("E26").Value = '  +0.13%  '
# Row 27
$ws.Range("E27").Value = '  +0.03%  '
# Row 28
$ws.Range("E28").Value = '  +1.32%  '
# Row 29
$ws.Range("E29").Value = '  +1.78%  '
# Row 30
$ws.Range("E30").Value = '  -0.76%  '
# Row 31
$ws.Range("D31").Value = '''0.0459'
$ws.Range("E31").Value = '  +33.55%  '
# Row 32
$ws.Range("D32").Value = '''50.96'
$ws.Range("E32").Value = '  +2.32%  '
# Row 33
$ws.Range("D33").Value = '''33.98'
$ws.Range("E33").Value = '  -3.01%  '
# Row 34
$ws.Range("E34").Value = '  +5.62%  '
# Row 35
$ws.Range("D35").Value = '''0.0832'
$ws.Range("E35").Value = '  +1.04%  '
# Row 36
$ws.Range("E36").Value = '  +0.01%  '
# Row 37
$ws.Range("E37").Value = '  -0.01%  '
# Row 38
$ws.Range("D38").Value = '''4.89'
$ws.Range("E38").Value = '  -2.27%  '
# Row 39
$ws.Range("E39").Value = '  -0.03%  '
# Row 40
$ws.Range("D40").Value = '''18.40'
$ws.Range("E40").Value = '  -3.85%  '
# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '''23.68'
$ws.Range("E41").Value = '  +2.98%  '
# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''2.58'
$ws.Range("E42").Value = '  +5.51%  '
# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''128.32'
$ws.Range("E43").Value = '  +1.06%  '
# Row 44
$ws.Range("E44").Value = '  +1.53%  '
# Row 45
$ws.Range("E45").Value = '  +0.51%  '
# Row 46
$ws.Range("E46").Value = '  +0.83%  '
# Row 47
$ws.Range("D47").Value = '''2.077.79'
$ws.Range("E47").Value = '  -0.60%  '
# Row 48
$ws.Range("E48").Value = '  +3.93%  '
# Row 49
$ws.Range("D49").Value = '''5.69'
$ws.Range("E49").Value = '  +3.05%  '
# Row 50
$ws.Range("D50").Value = '''0.945'
$ws.Range("E50").Value = '  +9.31%  '
# Row 51
$ws.Range("D51").Value = '''60.62'
$ws.Range("E51").Value = '  +1.34%  '
